$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.435.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.62%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.864.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4772"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3755"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07322"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9341"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07820"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.872.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.428"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.555"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.85%  "

$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008882"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.80%  "

$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.480.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.10%  "

$ws.Range("E22").Value = "  +1.00%  "

$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.943"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.69%  "

$ws.Range("E26").Value = "  +1.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.015"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.923"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08888"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.326"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.213"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7529"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.595"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.731"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.29%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02038"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.31%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.119"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05262"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.25%  "

$ws.Range("E39").Value = "  +0.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5315"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.057"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.576"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1522"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4798"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.95%  "

$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.655"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06078"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9182"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.40%  "
